$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newRows = @(
    @('JuggAdaptingArmor','Adaptive Armor','Takes 2 less damage from all attacks for 2 turns after taking damage.','None','Juggernaut'),
    @('JuggDamageControl','Damage Control','Grants 2 additional armor.','None','Juggernaut'),
    @('JuggDimensionalShift','Shields Overload','Expending all the shield energy within a short time, provides immunity to damage this turn, while in this state, the unit also fools the enemies as an easy target.','None','Juggernaut'),
    @('JuggDisorientationField','Disorientation Field','Disorients nearby enemies at the end of each turn.','None','Juggernaut'),
    @('JuggDistortionField','Group Protective Field','Produces a large protective field that reduces accuracy by 15% against any nearby squadmates.','None','Juggernaut'),
    @('JuggEnergyShield','Energy Shield','Generates a 5hp shield for nearby allies throughout the enemy turn.','None','Juggernaut'),
    @('JuggHeavyArmor','Heavy Armor','The armor absorbs damage taken, protecting the wearer from injury.','None','Juggernaut'),
    @('JuggHighCover','One for all','Become a heavy cover until an offensive action is taken.','Unknown','Juggernaut'),
    @('JuggProtect','Defender''s Will','Goes into a defensive stance, attacking any enemy who attacks in range until the soldier runs out of ammunition','Primary','Juggernaut'),
    @('JuggRegen','Feel No Pain','HP damage taken instantly gets regenerated as shields.','None','Juggernaut'),
    @('JuggResilient','Resilient','Become 15% harder to land a critical shot, and an additional 30% if flanked.','Unknown','Juggernaut'),
    @('JuggShockAbsorbantArmor','Shock Absorbant Armor','Reduce damage taken by 66% from nearby sources.','None','Juggernaut'),
    @('JuggThreatPrediction','Shaped Armor','Increases dodge by 75% for a turn','None','Juggernaut'),
    @('JackArmorPadding','Armor Padding','Gain +3 HP and receive an additional +10 Defense bonus. Less susceptible to stuns.','None','Trooper'),
    @('JackArtOfWar','Art of War','Deal +2 bonus damage to flanked enemies, with a +15% critical chance bonus.','None','Trooper'),
    @('JackDeathSentence','Death Sentence','Brace yourself. Fire a reaction shot for every enemy action within a radius. -20 Defense penalty. 5-turn cooldown.','Primary','Trooper'),
    @('JackLastStand','Survivor','Falling in battle heals the Trooper for 5HP and grants a +20 bonus to Will and Defense for the remaining battle. Once per battle.','Primary','Trooper'),
    @('JackLethalShot','Lethal Shot','Deal +1 bonus damage on attacks with a +5% Critical chance. Your attacks always graze enemies.','Primary','Trooper'),
    @('JackMarkTarget','Mark Target','Mark a target for 2 turns. All allies gain +15 Aim bonus when attacking the target.','Primary','Trooper'),
    @('JackPiercingBarrage','Piercing Barrage','Unleash a barrage of piercing rounds, damaging cover and enemies in a straight line. Penalty to aim and damage. Costs 2 ammo. 4-turn cooldown.','Primary','Trooper'),
    @('JackScatter','Spread Out','Nearby allies gain a 20% Dodge bonus for the next 2 turns, and a temporary boost to mobility. 4-turn cooldown.','Primary','Trooper'),
    @('JackSurvivalInstinct','Survival Instinct','Getting hit by an attack grants a +25 Defense and +5 Aim for one turn.','None','Trooper'),
    @('JackTactics','Aggressive Tactics','Increases target''s Aim by +5, Damage by +3, and a provides a great boost to mobility for one turn, with a -10 Defense penalty. 4 turn cooldown.','Primary','Trooper'),
    @('JackTrueShot','Lucky Shot','Unleash a focused shot from your pistol that ignores enemy armor. Guaranteed to hit. 3-turn cooldown.','Secondary','Trooper'),
    @('KeenEye','Entrench','Gain a +5 Defense bonus and a 5% bonus to critically hit enemies.','None','Trooper'),
    @('BlastShot','Share Blast','Fire a blast of energy in a narrow, long path in front of you. Can recieve +2 from HDD.','Primary','CPU'),
    @('HDD','HDD','Sacrifice 15 Defense for +2 damage, +20 aim and +25 crit. Can chain into other abilities! This health cannot be recovered until your next turn','None','CPU'),
    @('NepBull','Nep Bull','Self-use drink that restores up to 4 health to the user. 2 Uses per battle.','None','CPU'),
    @('NextGen','Next Gen','One of the CPU''s Gamechangers. Strike all foes in range with a pistol. Inflicts stun on FIRST TARGET. Can chain from HDD.','Secondary','CPU'),
    @('PartBreak','Part Break','Attack that focuses on breaking foe''s armor, shredding up to 6 points. 3 Uses.','Primary','CPU'),
    @('Patience','Patience','Wait for an opening from the enemy. This attack deals 2 bonus damage and is guaranteed to hit. Can chain from HDD for +4 Damage.','Primary','CPU'),
    @('Prioritize','Prioritize','Mark a taget, giveing everyone (Including yourself) +10 Aim on it.','None','CPU'),
    @('Sentient','Sentient','Ends turn with +10 Armor until the start of the next turn. Helps avoid death in tight spots.','None','CPU'),
    @('Smite','Smite','Pistol attack that causes the enemy to panic. Deals no damage however (Cannot gain +2 from HDD)','Primary','CPU'),
    @('Swiftness','Swiftness','Increases Mobility stat by 2 points.','None','CPU'),
    @('Teleport','Warp','One of the CPU''s Gamechangers. Allows the user to warp to another location.','None','CPU'),
    @('NecroCombustion','Combustion','Cause a creature''s blood to boil. Does 2-4 guaranteed psi damage and sets them on fire.','Secondary','Necromancer'),
    @('NecroDeathTouch','Death Touch','The Necromancer has gained a superhuman understanding of life and death. Touch a living creature to grab its soul and pull it out of its body. It dies instantly. It can be used once per battle, from melee range only.','Secondary','Necromancer'),
    @('NecroDrainLife','Drain Life','Shoots a target enemy and transfers half of the damage done back to the Necromancer as health.','Primary','Necromancer'),
    @('NecroEldritchStorm','Eldritch Storm','The necromancer has become an avatar of death. Focus all your power into one spell that will Set on Fire, Poison,and cause Acid Burn to all enemies in a large radius. Does no direct damage.','Secondary','Necromancer'),
    @('NecroFireballAbility','Fireball','Launch a Fireball, dealing explosive damage to everyone in the blast radius and setting everything on fire.','Primary','Necromancer'),
    @('NecroHallucination','Hallucinations curse','The necromancer curses a creature, causing them to see their own team as enemies, and then overloads its mind with impending doom. Surrounded from all sides, they are guaranteed to panic, unless they have full immunity to panic. There is a good chance they will shoot their own team.','None','Necromancer'),
    @('NecroLightning','Lightning','Shoots a lightning bolt at a target location. Damages everyone caught in the radius, stuns organic units and causes robotic ones to reboot. Costs 2 action points.','Primary','Necromancer'),
    @('NecroMassRaiseDead','Army of the dead','Raises all dead humans in a large radius as zombies.','Secondary','Necromancer'),
    @('NecroPoison','Poison','Poison any living, organic creature. Does 1 damage per turn, -30 to aim and -4 to move. The spell doesn''t end your turn and has no cooldown.','Secondary','Necromancer'),
    @('NecroRaiseDead','Raise Zombie','Raise a dead human or Advent as a zombie with 8 HP.','None','Necromancer'),
    @('NecroRestoration','Restoration','Regenerate 1 HP per turn, for 3 turns, to any ally. \n Shamans are said to treat ailments/illness by mending the soul. Alleviating traumas affecting the soul/spirit restores the physical body of the individual to balance and wholeness.','None','Necromancer'),
    @('NecroSoulSpear','Soul Spear','Sorcery that fires souls shaped into a spear. Pierces enemies, causing heavy damage. This spell was said to have been devised by a master sorcerer, but his name is long forgotten.','Primary','Necromancer'),
    @('NecroSpiritArmorAbility','Bone Armor','Adds 3 armor to a friendly target. This bonus lasts until the end of the mission and can be stacked. 4 turns cool down.','None','Necromancer'),
    @('NecroStrenghtenUndead','Undead Link','Enhance the link between the Necromancer and one of his zombies, significantly increasing all of its life stats.','None','Necromancer'),
    @('AbsorptionField','Channeling Field','Every time you are targeted, part of the energy is channeled to your <Ability:WeaponName>. This energy is unleashed with your next Standard Shot.','None','Shen''s Last Gift'),
    @('AdaptiveAim','Adaptive Aim','When Overdrive is active, Standard Shots do not incur recoil penalties.','Secondary','Shen''s Last Gift'),
    @('Arsenal','Arsenal','The BIT can equip and fire heavy weapons.','Unknown','Shen''s Last Gift'),
    @('Bombard','Bombard','Launch the BIT to any visible location where it releases a powerful explosive blast. Advanced BITs do more damage.','Secondary','Shen''s Last Gift'),
    @('Bulwark','Bulwark','Gain a bonus point of Armor, and always provide high cover to all adjacent squadmates.','None','Shen''s Last Gift'),
    @('CapacitorDischarge_Shen','Capacitor Discharge','Send the GREMLIN to a location where it emits a substantial electric discharge, damaging and potentially stunning all nearby units. Robotic units take more damage.','Secondary','Shen''s Last Gift'),
    @('CombatProtocol_Shen','Combat Protocol','Send the GREMLIN to an enemy to jolt them, dealing guaranteed damage, which is increased against robotic enemies.','Secondary','Shen''s Last Gift'),
    @('HunterProtocol','Hunter Protocol','When an enemy is revealed, there is a <Ability:HunterProtocol_Chance/>% chance to take a free Overwatch shot.','Primary','Shen''s Last Gift'),
    @('Intimidate','Intimidate','When targeted by an attack, the enemy has a chance to panic.','None','Shen''s Last Gift'),
    @('MedicalProtocol_Shen','Medical Protocol','The GREMLIN can perform healing actions remotely: GREMLIN Heal and GREMLIN Stabilize. The GREMLIN has four charges. If a medikit is equipped, the GREMLIN will gain an additional charge.','Secondary','Shen''s Last Gift'),
    @('Nova','Nova','Release a blast of energy, damaging all nearby enemies. This attack has no cost or cooldown, but consecutive Novas will damage the SPARK.','Unknown','Shen''s Last Gift'),
    @('Overdrive','Overdrive','Take three actions this turn, and no action is turn-ending. Multiple Standard Shots incur a small recoil penalty.','Primary','Shen''s Last Gift'),
    @('Rainmaker','Rainmaker','Equipped heavy weapons deal +<Ability:Rainmaker_BonusDamage/> damage and have increased area of effect.','Unknown','Shen''s Last Gift'),
    @('Repair','Repair','Send the BIT to a damaged robotic ally and repair it. Advanced BITs can repair slightly more.','Secondary','Shen''s Last Gift'),
    @('Sacrifice','Sacrifice','Generate a protective field which redirects any attacks against allies inside it towards you. Gain bonus defense and armor while active.','Unknown','Shen''s Last Gift'),
    @('Strike','Strike','Attack any enemy within movement range with a powerful melee attack.','Unknown','Shen''s Last Gift'),
    @('WreckingBall','Wrecking Ball','While Overdrive is active, break through walls and cover when moving.','None','Shen''s Last Gift')
)

$startRow = 488
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $rowData = $newRows[$i]
    $targetRow = $startRow + $i
    for ($j = 0; $j -lt $rowData.Count; $j++) {
        $ws.Cells.Item($targetRow, $j + 1).Value = $rowData[$j]
    }
}
